# Testexercise5.xlsx - "Final version including 20 Test cases"
#
# Fixes a couple of wording issues on the "calculate_price_test" sheet
# (typo "Boundari" -> "Boundary", and a clearer description for the
# float-total test case), and leaves the sheet's active-cell selection
# on C26 (it was previously left on D26).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("calculate_price_test")

# Typo fix: "Boundari" -> "Boundary" (appears in the Title column of the
# test-case rows built around the Total=0 boundary value test). Row 16 is
# fixed first, then row 2 (matches the order the shared-string table ends
# up in on save).
$ws.Range("B16").Value = "Boundary value test, Total=0"
$ws.Range("C16").Value = "Boundary value test, Total=0"
$ws.Range("B2").Value  = "During easter holiday, first tier, Boundary value test"

# Reworded objective for the "float total" test case.
$ws.Range("C21").Value = "The argument total has a type of float."

# Move the sheet's remembered selection from D26 to C26.
$ws.Activate()
$ws.Range("C26").Select()
